$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.12%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '39.72'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.55%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.153'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.71%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08104'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.66%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.946'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-5.03%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.160'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.30%'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9294'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.02%'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1399'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.47%'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1921'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.23%'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09151'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.55%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03494'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.49%'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09826'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.13%'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001392'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.89%'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005848'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.77%'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.944'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.05%'

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.228'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.28%'

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.92%'

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.75%'

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.45%'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.683'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.37%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2423'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-7.43%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04379'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.96%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001233'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.71%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004353'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '4.31%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.06%'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004003'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-10.00%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02034'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-3.75%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05058'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.10%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007388'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.07%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009762'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-3.59%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1363'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.40%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002130'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.06%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009501'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.75%'

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006346'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.57%'

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.05%'

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002736'

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-18.69%'

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.05%'

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.05%'
